$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Assign work to Gabo (row 4): set the "Team Member" column (B) value to "EDT"
$ws.Range("B4").Value = "EDT"

# Move the active selection to C4
$ws.Range("C4").Select()
